# Ran code for averaged intensities on spiral schemes
#
# The averaging run now also covers the Gaussian-Quadrature scheme and three
# new spiral schemes. Those land right after the "Ring Perpendicular to TD"
# row, which pushes the pre-existing NoRotation-tilt60deg / Rotation-* /
# HexGrid-* rows down the table; the three HexGrid rows are re-appended at
# the bottom (rows 17-19) to keep the same information.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with three more data rows, copying the formatting
# (bold/centered/bordered column A) from an existing data row.
$ws.Range("A16:M16").Copy() | Out-Null
$ws.Range("A17:M19").PasteSpecial(-4122) | Out-Null

$labels = @(
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
